$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows before the old last row (row 65), pushing it down to row 69.
$ws.Rows("65:68").Insert()

# Match the formatting (border/alignment style) used by column A on the
# surrounding data rows for the newly inserted rows.
$ws.Range("A64").Copy()
$ws.Range("A65:A68").PasteSpecial(-4122)

# Final state for rows 34..69 (columns A-E).
# A = sequential counter (row-2), B = 15000001 (constant), C = static-entry id,
# D = "show" (constant), E = terminology code "ITCHxxxxx".
$finalRows = @(
    @(34, 32, 14000034, "ITCH00001"),
    @(35, 33, 14000035, "ITCH00002"),
    @(36, 34, 14000036, "ITCH00003"),
    @(37, 35, 14000038, "ITCH00004"),
    @(38, 36, 14000039, "ITCH00005"),
    @(39, 37, 14000040, "ITCH00006"),
    @(40, 38, 14000041, "ITCH00007"),
    @(41, 39, 14000042, "ITCH00009"),
    @(42, 40, 14000044, "ITCH00010"),
    @(43, 41, 14000046, "ITCH00011"),
    @(44, 42, 14000047, "ITCH00012"),
    @(45, 43, 14000049, "ITCH00013"),
    @(46, 44, 14000050, "ITCH00014"),
    @(47, 45, 14000051, "ITCH00015"),
    @(48, 46, 14000052, "ITCH00016"),
    @(49, 47, 14000053, "ITCH00017"),
    @(50, 48, 14000054, "ITCH00018"),
    @(51, 49, 14000055, "ITCH00019"),
    @(52, 50, 14000056, "ITCH00020"),
    @(53, 51, 14000057, "ITCH00021"),
    @(54, 52, 14000058, "ITCH00022"),
    @(55, 53, 14000059, "ITCH00023"),
    @(56, 54, 14000060, "ITCH00024"),
    @(57, 55, 14000061, "ITCH00025"),
    @(58, 56, 14000063, "ITCH00026"),
    @(59, 57, 14000064, "ITCH00027"),
    @(60, 58, 14000066, "ITCH00028"),
    @(61, 59, 14000068, "ITCH00029"),
    @(62, 60, 14000069, "ITCH00030"),
    @(63, 61, 14000070, "ITCH00031"),
    @(64, 62, 14000071, "ITCH00032"),
    @(65, 63, 14000072, "ITCH00033"),
    @(66, 64, 14000073, "ITCH00034"),
    @(67, 65, 14000074, "ITCH00035"),
    @(68, 66, 14000075, "ITCH00036"),
    @(69, 67, 14000076, "ITCH00040")
)

foreach ($r in $finalRows) {
    $rowNum = $r[0]
    $aVal = $r[1]
    $cVal = $r[2]
    $eVal = $r[3]

    $ws.Cells.Item($rowNum, 1).Value = $aVal
    $ws.Cells.Item($rowNum, 2).Value = 15000001
    $ws.Cells.Item($rowNum, 3).Value = $cVal
    $ws.Cells.Item($rowNum, 4).Value = "show"
    $ws.Cells.Item($rowNum, 5).Value = $eVal
}
